# Update the "last_edited_time" column (D) for the Notion export rows.
# The underlying shared-string table was regenerated by the source tool, which
# (after de-duplication) resolves to these per-row ranges of new timestamps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D4").Value = "2024-08-31T05:43:00.000Z"
$ws.Range("D5:D6").Value = "2024-08-31T05:39:00.000Z"
$ws.Range("D7:D57").Value = "2024-08-31T05:40:00.000Z"
$ws.Range("D58:D102").Value = "2024-08-31T05:41:00.000Z"
$ws.Range("D103:D112").Value = "2024-08-31T05:42:00.000Z"
